$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.25385405913030923,0.23493005682721724,0.23620012671432528,0.27954658841439545,0.2537648004717664,0.24694905260139216,0.24403045852383065,0.24453365254631745,0.24666798958238953,0.24583811340227155,0.25179985559160456,0.27242573395598507,0.25602693100864793,0.2646881981053287,0.2489368314622166,0.24561692903668075,0.26116868065074955,0.2500092037952017,0.24614755975586275,0.23129864235600842,0.2321948147156281,0.2394840352627732,0.23211265813536103,0.22104525806303954,0.20528601275261193,0.19613881844534028,0.19242049848292878,0.18899065532511425,0.19445511934915635,0.20080732820829467,0.1717461099846498,0.1593696224863909,0.1441111575314644,0.14634845049670148,0.16944285308808746,0.1841518431774174,0.17593322036020123,0.15409054491860893,0.17470571549463848,0.1731983141123689,0.16968407162092386,0.1873128974756765,0.19751667771031614,0.20540850833403387,0.20945865620530651,0.22027734456284323,0.2160829303290156,0.21019801485231562,0.19735210741154838,0.21185634372326648,0.22786783984417752,0.21648475786800309,0.2419200744799142,0.236663641795989,0.2455847376168069,0.23444219756433507,0.22259836365105878,0.2500403933274996,0.25190069002649057,0.25347578651007163,0.2641012279369823)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $values[$i]
}
